# Pre-App Recommend Condition Letter
# Commit: "Paragraph number and year changes, plus additional changes"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert three additional blank "KeyHeadDetails" paragraphs between
#    "Our Ref: <Primary Reference Number>" and the paragraph holding the
#    four tab characters.
# ---------------------------------------------------------------------
$rng = $d.Content
$foundRef = $rng.Find.Execute("Our Ref:", $true, $false, $false, $false, $false, $false, 1, $false, "", 0)
if ($foundRef) {
    $ourRefPara = $rng.Paragraphs(1)
    $ourRefPara.Range.InsertParagraphAfter()
    $ourRefPara.Range.InsertParagraphAfter()
    $ourRefPara.Range.InsertParagraphAfter()
}

# ---------------------------------------------------------------------
# 2) "NATIONAL PLANNING POLICY FRAMEWORK 2021" -> "...2023"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("NATIONAL PLANNING POLICY FRAMEWORK 2021", $true, $false, $false, $false, $false, $false, 1, $false, "NATIONAL PLANNING POLICY FRAMEWORK 2023", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) <Proposal Description> paragraph becomes italic (in addition to bold)
# ---------------------------------------------------------------------
$rng2 = $d.Content
$foundPD = $rng2.Find.Execute("<Proposal Description>", $true, $false, $false, $false, $false, $false, 1, $false, "", 0)
if ($foundPD) {
    $pdPara = $rng2.Paragraphs(1)
    $pdPara.Range.Font.Italic = 1
    $pdPara.Range.Font.ItalicBi = 1
}

# ---------------------------------------------------------------------
# 4) "Recommend Archaeological Condition(s)" -> "Recommend Archaeological Condition"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Recommend Archaeological Condition(s)", $true, $false, $false, $false, $false, $false, 1, $false, "Recommend Archaeological Condition", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) NPPF paragraph number changes
# ---------------------------------------------------------------------
$d.Content.Find.Execute("NPPF paragraph 194 says", $true, $false, $false, $false, $false, $false, 1, $false, "NPPF paragraph 200 says", 2) | Out-Null
$d.Content.Find.Execute("NPPF paragraphs 190 and 197", $true, $false, $false, $false, $false, $false, 1, $false, "NPPF paragraphs 195 and 203", 2) | Out-Null
$d.Content.Find.Execute("paragraph 205 of the NPPF", $true, $false, $false, $false, $false, $false, 1, $false, "paragraph 211 of the NPPF", 2) | Out-Null

# ---------------------------------------------------------------------
# 6) Second <Casework Officer> (below "Yours sincerely") becomes bold
# ---------------------------------------------------------------------
$ySincere = $d.Content
$foundYS = $ySincere.Find.Execute("Yours sincerely", $true, $false, $false, $false, $false, $false, 1, $false, "", 0)
if ($foundYS) {
    $afterSincerely = $d.Range($ySincere.End, $d.Content.End)
    $foundCO = $afterSincerely.Find.Execute("<Casework Officer>", $true, $false, $false, $false, $false, $false, 1, $false, "", 0)
    if ($foundCO) {
        $coPara = $afterSincerely.Paragraphs(1)
        $coPara.Range.Font.Bold = 1
        $coPara.Range.Font.BoldBi = 1
    }
}
